# Refitting NCDEs to individual patients (for manuscript figure)
#
# Adds a new "Label" column (H) that flags each patient row as Control (0)
# or MDD (1), and refreshes several prediction/error/cross-entropy-loss
# values from the re-fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header -----------------------------------------------
$ws.Range("H1").Value = "Label"

# --- Block 1: 100 iterations (rows 2-11) -------------------------------
# Controls (46, 28, 13, 50, 51) -> Label = 0
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
# MDD (27, 47, 13, 25, 5) -> Label = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# --- Block 2: 200 iterations (rows 12-21) ------------------------------
# Controls (46, 28, 13, 50, 51) -> Label = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
# MDD (27, 47, 13, 25, 5) -> Label = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1

# --- Refreshed Prediction / Error / Cross Entropy Loss values ----------
$ws.Range("D2").Value = 0.4368654151023003
$ws.Range("E2").Value = 0.4368654151023003

$ws.Range("D4").Value = 0.4109964014357663
$ws.Range("E4").Value = 0.4109964014357663

$ws.Range("D5").Value = 0.4827490253937415
$ws.Range("E5").Value = 0.4827490253937415

$ws.Range("D6").Value = 0.7778153725582956
$ws.Range("E6").Value = 0.7778153725582956

$ws.Range("D8").Value = 0.5556080447587548
$ws.Range("E8").Value = 0.4443919552412452

$ws.Range("D9").Value = 0.4706064302781766
$ws.Range("E9").Value = 0.5293935697218234

$ws.Range("D10").Value = 0.6810187797930835
$ws.Range("E10").Value = 0.3189812202069165

$ws.Range("F11").Value = 0.6924229264259338
